# Updates the cryptos price/volume snapshot (GitHub Actions scheduled refresh).
# Columns: A=rank index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Values that parse as plain decimals (e.g. "213.97") are entered with a
# leading apostrophe, mirroring how Excel keeps user-typed numeric-looking
# text as text instead of auto-converting it to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.969.87'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.633.19'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''213.97'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').Value = '''0.504'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').Value = '''18.47'
$ws.Range('E10').Value = '  -6.06%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '1.859.94'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.19'
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.618.74'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('D17').Value = '25.979.67'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '''61.61'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '''190.42'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').Value = '''4.24'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').Value = '''9.56'
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').Value = '''6.13'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '''143.10'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -3.31%  '
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').Value = '''15.22'
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('E31').Value = '  -3.26%  '
$ws.Range('E32').Value = '  -3.08%  '
$ws.Range('E33').Value = '  -4.14%  '
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('D36').Value = '1.134.37'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '''0.865'
$ws.Range('E37').Value = '  -4.58%  '
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('E39').Value = '  -3.24%  '
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('D41').Value = '''98.60'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('E43').Value = '  -5.18%  '
$ws.Range('D44').Value = '1.769.97'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('D46').Value = '''55.05'
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = '''7.56'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('E51').Value = '  +0.03%  '
